$wb = $excel.ActiveWorkbook

$wsSite      = $wb.Worksheets.Item("site")
$wsOuting    = $wb.Worksheets.Item("outing")
$wsCapture   = $wb.Worksheets.Item("capture")
$wsRecapture = $wb.Worksheets.Item("recapture")

# --- outing sheet: move selection, clear two stray checkbox values ---
$wsOuting.Activate()
$wsOuting.Range("I5").ClearContents() | Out-Null
$wsOuting.Range("F6").ClearContents() | Out-Null
$wsOuting.Range("I5").Select() | Out-Null

# --- capture sheet: move selection, add a new checked box at K5 ---
$wsCapture.Activate()
$wsCapture.Range("K5").Value2 = $true
$wsCapture.Range("B6").Select() | Out-Null

# --- recapture sheet: move selection, clear some checkbox values, becomes active tab ---
$wsRecapture.Activate()
$wsRecapture.Range("B5:D5").ClearContents() | Out-Null
$wsRecapture.Range("G6").ClearContents() | Out-Null
$wsRecapture.Range("E8").Select() | Out-Null

# Final active sheet is "recapture" (activeTab=3, tabSelected on its sheetView)
$wsRecapture.Activate()
